$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "...Oracle, Git." -> "...Oracle, Git, Microsoft Office."
# (Technical skills paragraph) - insert ", Microsoft Office" right after
# "Git" (before the trailing period), then append the period back as its
# own run so the three pieces of text end up in separate runs, just like
# a user placing the cursor and typing the addition in Word.
# ---------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Oracle, Git.") | Out-Null
$insertPos = $findRng.Start + ("Oracle, Git".Length)

$insRng = $d.Range($insertPos, $insertPos)
$insRng.InsertAfter(", Microsoft Office")
# Nudge formatting (no visual change) so the new text is materialised as
# its own run(s) instead of being silently folded back into its neighbour.
$insRng.Font.Bold = $true
$insRng.Font.Bold = $false

# ---------------------------------------------------------------------
# Change 2: drop the forced "space after = 0" direct formatting on the
# "Foodomnia" job-heading paragraph so it falls back to the style's
# default spacing (8pt / 160 twips, i.e. Normal style's w:spacing w:after).
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Foodomnia*") {
        $p.Format.SpaceAfter = 8
        break
    }
}

# ---------------------------------------------------------------------
# Change 3: remove the whole bullet paragraph about automating the
# weekly financial/operational report process.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Automated the weekly process*") {
        $p.Range.Delete()
        break
    }
}
